$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Replace the old combined sentence with the new combined sentence
# (single Find/Replace call keeps this as one run for now; we will re-split
# it into the individual runs shown by the diff afterwards).
# ---------------------------------------------------------------------------
$oldSentence = "Although I found this work very challenging, I believe it helped reinforce some basic java fundamentals as well as key concept we learn in class such as how to read from a file and throwing exceptions."
$newSentence = "Additionally, I also ran into roadblocks where certain Start Destinations and End Destinations that had flights connecting them would return an error or not return a path whiles others would. For instance, Accra, Ghana to Winnipeg, Canada would result in no route although there is one but for another case like Accra, Ghana to London, United Kingdom it would give me a solution. Although I found this work very challenging, I believe it helped reinforce some basic java fundamentals as well as key concept we learn in class such as how to read from a file and throwing exceptions."

$r1 = $d.Content
$r1.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: Split the paragraph into two paragraphs right before the second
# "Although I found this work very challenging" sentence.
# ---------------------------------------------------------------------------
$rSplit = $d.Content
$rSplit.Find.Execute("Although I found this work very challenging") | Out-Null
$splitPoint = $d.Range($rSplit.Start, $rSplit.Start)
$splitPoint.InsertParagraphAfter()

Write-Host "After split, paragraph count:" $d.Paragraphs.Count
